$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text so values like "1.00" are not
# coerced to numbers (source data is literal text, not numeric).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "70.279.29"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "3.790.31"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "609.02"
$ws.Range("E5").Value = "  +3.82%  "
$ws.Range("D6").Value = "188.02"
$ws.Range("E6").Value = "  +14.89%  "
$ws.Range("D7").Value = "3.784.27"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "0.643"
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "0.738"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("D12").Value = "57.84"
$ws.Range("E12").Value = "  +11.21%  "
$ws.Range("D13").Value = "0.0000301"
$ws.Range("E13").Value = "  -4.64%  "
$ws.Range("D14").Value = "10.89"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "4.382.71"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "3.786.92"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "19.81"
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").Value = "13.20"
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "70.161.79"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "423.42"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "90.99"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").Value = "3.11"
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "11.47"
$ws.Range("E26").Value = "  +3.98%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "13.18"
$ws.Range("E27").Value = "  -3.58%  "
$ws.Range("D28").Value = "4.10"
$ws.Range("E28").Value = "  +4.66%  "
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("D30").Value = "9.76"
$ws.Range("E30").Value = "  -5.28%  "
$ws.Range("D31").Value = "33.61"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("D32").Value = "7.64"
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").Value = "12.74"
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").Value = "45.09"
$ws.Range("E35").Value = "  -5.32%  "
$ws.Range("D36").Value = "622.17"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").Value = "65.33"
$ws.Range("E37").Value = "  -4.58%  "
$ws.Range("D38").Value = "0.0₃0930"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").Value = "0.414"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "3.13"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("D45").Value = "3.08"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").Value = "0.0453"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "9.55"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.137"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "3.26"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.828.56"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("E51").Value = "  -1.80%  "
